# The workbook has many "blank" numeric cells (formatted/typed as numbers
# but holding no value). This commit fills each of those blank cells with
# an explicit 0, without touching any cell that already has content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "F4","H4","K4","M4","P4","R4",
    "P5","R5",
    "D6","I6","J6","O6","S6","T6",
    "D8","I8","J8","O8","S8","T8",
    "F9","H9","K9","M9","P9","R9",
    "F10","H10","K10","M10","P10","R10",
    "D11","I11","J11","O11","S11","T11",
    "D13","I13","J13","O13","S13","T13",
    "F14","H14","K14","M14","P14","R14",
    "F15","H15","K15","M15","P15","R15",
    "D16","I16","J16","O16","S16","T16",
    "D18","I18","J18","O18","S18","T18",
    "P19","R19",
    "F20","H20","K20","M20","P20","R20"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 0
}
